$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.98"
$ws.Range("E2").Value = "'0.87%"
$ws.Range("D3").Value = "'29.19"
$ws.Range("E3").Value = "'-2.84%"
$ws.Range("D4").Value = "'5.146"
$ws.Range("E4").Value = "'-0.03%"
$ws.Range("D5").Value = "'0.05791"
$ws.Range("E5").Value = "'2.01%"
$ws.Range("D6").Value = "'6.644"
$ws.Range("E6").Value = "'1.73%"
$ws.Range("D8").Value = "'0.8567"
$ws.Range("E8").Value = "'1.56%"
$ws.Range("D9").Value = "'0.8641"
$ws.Range("E9").Value = "'0.91%"
$ws.Range("D10").Value = "'0.1373"
$ws.Range("E10").Value = "'2.99%"
$ws.Range("D11").Value = "'0.07074"
$ws.Range("E11").Value = "'2.37%"
$ws.Range("D12").Value = "'0.03292"
$ws.Range("E12").Value = "'14.08%"
$ws.Range("D13").Value = "'0.09367"
$ws.Range("E13").Value = "'-0.13%"
$ws.Range("D14").Value = "'0.001528"
$ws.Range("E14").Value = "'0.74%"
$ws.Range("B15").Value = "TigerCash"
$ws.Range("C15").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D15").Value = "'0.006069"
$ws.Range("E15").Value = "'-1.17%"
$ws.Range("B16").Value = "LEO"
$ws.Range("C16").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D16").Value = "'3.487"
$ws.Range("E16").Value = "'-0.60%"
$ws.Range("B17").Value = "BTSEToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D17").Value = "'2.182"
$ws.Range("E17").Value = "'2.25%"
$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D18").Value = "'0.0006023"
$ws.Range("E18").Value = "'-0.01%"
$ws.Range("E19").Value = "'1.59%"
$ws.Range("D20").Value = "'0.03345"
$ws.Range("E20").Value = "'1.92%"
$ws.Range("D21").Value = "'0.1282"
$ws.Range("E21").Value = "'-1.62%"
$ws.Range("D22").Value = "'3.180"
$ws.Range("E22").Value = "'-11.89%"
$ws.Range("D23").Value = "'0.04143"
$ws.Range("E23").Value = "'-0.59%"
$ws.Range("E25").Value = "'1.07%"
$ws.Range("E27").Value = "'2.51%"
$ws.Range("E28").Value = "'3.37%"
$ws.Range("D40").Value = "'0.03731"
$ws.Range("E40").Value = "'0.45%"
$ws.Range("D41").Value = "'0.005803"
$ws.Range("E41").Value = "'70.77%"
$ws.Range("D42").Value = "'0.1070"
$ws.Range("E42").Value = "'1.03%"
$ws.Range("D43").Value = "'0.002197"
$ws.Range("E43").Value = "'-3.92%"
$ws.Range("D44").Value = "'0.009167"
$ws.Range("E44").Value = "'-6.84%"
$ws.Range("D45").Value = "'0.00005272"
$ws.Range("E45").Value = "'3.37%"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E47").Value = "'-42.01%"
$ws.Range("E48").Value = "'-22.63%"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("E50").Value = "'-0.01%"